# CSCohortTest.xlsx update
# "concept of permanency of click added" -- extends the Year1 cohort-tracking
# table with students 29-45 (rows 31-47), following the same repeating
# 10-row pattern already present in the sheet, and records a newly observed
# course code (E64000) for student #5's (row 7) later years as well as for
# the newly added student #33 (row 35).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- 1. Update row 7 (student 5): F7:K7 change from E28600 -> E64000 -------
$ws.Range("F7:K7").Value2 = "E64000"

# --- 2. Append rows 31-47 (students 29-45), continuing the existing cycle --
$newRowsData = @(
    @(29,"MALE","YES","E28400","E28400","E28400","E28400","E28400","E28400","E28400","E28400"),
    @(30,"FEMALE","YES","E28400","E28400","E28400","E28400","E28600","E28600","E28600","E28600"),
    @(31,"MALE","YES","E28400","E10600","E10600","E10600","DropOut","DropOut","DropOut","DropOut"),
    @(32,"FEMALE","YES","E28400","E10600","E10600","E10600","E10600","E10600","E10600","E10600"),
    @(33,"MALE","YES","E28400","E28600","E64000","E64000","E64000","E64000","E64000","E64000"),
    @(34,"FEMALE","NO","E28400","E28600","E28600","E28600","E28600","Dismissed","E28600","E28600"),
    @(35,"MALE","NO","E28400","E28600","E28600","E28600","E28600","E28600","E28600","E28600"),
    @(36,"FEMALE","NO","E28400","E28600","E28600","E28600","E28600","E10600","E10600","E10600"),
    @(37,"MALE","NO","E28400","DropOut","DropOut","DropOut","DropOut","DropOut","DropOut","DropOut"),
    @(38,"FEMALE","NO","E28400","Dismissed","Dismissed","Dismissed","Dismissed","Dismissed","Dismissed","Dismissed"),
    @(39,"MALE","YES","E28400","E28400","E28400","E28400","E28400","E28400","E28400","E28400"),
    @(40,"FEMALE","YES","E28400","E28400","E28400","E28400","E28600","E28600","E28600","E28600"),
    @(41,"MALE","YES","E28400","E10600","E10600","E10600","DropOut","DropOut","DropOut","DropOut"),
    @(42,"FEMALE","YES","E28400","E10600","E10600","E10600","E10600","E10600","E10600","E10600"),
    @(43,"MALE","YES","E28400","E28600","E28600","E28600","E28600","E28600","E28600","E28600"),
    @(44,"FEMALE","NO","E28400","E28600","E28600","E28600","E28600","Dismissed","E28600","E28600"),
    @(45,"MALE","NO","E28400","E28600","E28600","E28600","E28600","E28600","E28600","E28600")
)

$startRow = 31
for ($i = 0; $i -lt $newRowsData.Length; $i++) {
    $rowNum = $startRow + $i
    $rowValues = $newRowsData[$i]
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $ws.Cells.Item($rowNum, $c + 1).Value2 = $rowValues[$c]
    }
}

# --- 3. Update the view: scroll position and active selection --------------
[void]$ws.Range("F113").Select()
$excel.ActiveWindow.ScrollRow = 96
$excel.ActiveWindow.ScrollColumn = 1
